# Update cryptos list values (price/volume) as scraped on the latest run,
# including a row swap between THORChain and ApeXProtocol (rows 48/49).
#
# NOTE: several "Price" values (column D) are plain decimal numbers (e.g.
# "586.92"). Excel's COM Range.Value setter auto-converts such text into a
# real number, which would change the cell's stored type away from the
# text type used throughout this sheet. To keep those values as text we
# briefly force a text number format while assigning them, then restore
# the original "General" format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.350.05"
$ws.Range("E2").Value = "  -4.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.563.25"
$ws.Range("E3").Value = "  -4.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "586.92"
$ws.Range("E5").Value = "  -4.68%  "

# Row 6 - Solana
Set-TextValue "D6" "185.04"
$ws.Range("E6").Value = "  -1.14%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.555.86"
$ws.Range("E7").Value = "  -4.80%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -4.27%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.11%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.670"
$ws.Range("E10").Value = "  -7.32%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -10.01%  "

# Row 12 - Avalanche
Set-TextValue "D12" "53.06"
$ws.Range("E12").Value = "  -6.84%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -12.31%  "

# Row 14 - Polkadot
Set-TextValue "D14" "9.79"
$ws.Range("E14").Value = "  -8.67%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.135.25"
$ws.Range("E15").Value = "  -4.63%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.564.83"
$ws.Range("E16").Value = "  -4.80%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.96%  "

# Row 18 - Chainlink
Set-TextValue "D18" "18.29"
$ws.Range("E18").Value = "  -5.70%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -7.04%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "66.224.42"
$ws.Range("E20").Value = "  -4.15%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  -7.51%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "395.45"
$ws.Range("E22").Value = "  -4.55%  "

# Row 23 - PancakeSwap
$ws.Range("E23").Value = "  -6.64%  "

# Row 24 - Litecoin
Set-TextValue "D24" "85.70"
$ws.Range("E24").Value = "  -4.24%  "

# Row 25 - RenderToken
Set-TextValue "D25" "11.23"
$ws.Range("E25").Value = "  +2.06%  "

# Row 26 - ImmutableX
Set-TextValue "D26" "2.89"
$ws.Range("E26").Value = "  -5.88%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.42"
$ws.Range("E27").Value = "  -3.57%  "

# Row 28 - LEO
Set-TextValue "D28" "6.04"
$ws.Range("E28").Value = "  -0.58%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -7.31%  "

# Row 30 - Filecoin
Set-TextValue "D30" "8.92"
$ws.Range("E30").Value = "  -8.29%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "31.07"
$ws.Range("E31").Value = "  -6.71%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "7.01"
$ws.Range("E32").Value = "  -4.89%  "

# Row 33 - Bittensor
Set-TextValue "D33" "623.71"
$ws.Range("E33").Value = "  +1.36%  "

# Row 34 - Cosmos
Set-TextValue "D34" "12.14"
$ws.Range("E34").Value = "  -4.91%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  -8.27%  "

# Row 36 - OKB
Set-TextValue "D36" "62.98"
$ws.Range("E36").Value = "  -5.14%  "

# Row 37 - InjectiveProtocol
Set-TextValue "D37" "41.17"
$ws.Range("E37").Value = "  -7.90%  "

# Row 38 - Dai
$ws.Range("E38").Value = "  +0.07%  "

# Row 39 - TheGraph
Set-TextValue "D39" "0.398"
$ws.Range("E39").Value = "  -2.76%  "

# Row 40 - PEPE
$ws.Range("D40").Value = "0.0₃0761"
$ws.Range("E40").Value = "  -12.14%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.29%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -8.05%  "

# Row 43 - Maker
$ws.Range("D43").Value = "3.011.73"
$ws.Range("E43").Value = "  +5.61%  "

# Row 44 - ThetaToken
$ws.Range("E44").Value = "  -9.03%  "

# Row 45 - Fetch.AI
Set-TextValue "D45" "2.52"
$ws.Range("E45").Value = "  -4.90%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -8.39%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  -8.26%  "

# Row 48/49 swapped: ApeXProtocol now ranks above THORChain
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D48" "3.04"
$ws.Range("E48").Value = "  -2.10%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D49" "8.55"
$ws.Range("E49").Value = "  -7.56%  "

# Row 50 - Monero
Set-TextValue "D50" "138.02"
$ws.Range("E50").Value = "  -2.89%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  -1.66%  "
